# Update the "Förändrad" (Changed) date column for rows 2-7 from
# 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 45212
